$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A1").Value = 45432
$ws.Range("D30").Value = 202.54
$ws.Range("D31").Value = 230
